# Append "Trade #17" row (sheet row 18) to both the "All Trades" and
# "base_strategy" worksheets, closing trade row r=18:
#   A=17 (n), B="2026-02-16" (text), C="22:59:02" (text),
#   D="base_strategy" (text), E="DOWN" (text), F=0.5 (n), G="" (blank),
#   H="OPEN" (text), I=0 (n), J=0 (n), K=100 (n), L=0 (n), M=0 (n),
#   N=0.6 (n), O="Normal spread capture: 19600 bps" (text), P="" (blank),
#   Q=0 (n)

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 18

    $ws.Cells.Item($row, 1).Value = 17

    # Column B holds a text value that looks like a date ("2026-02-16").
    # Pre-format the cell as Text so Excel does not auto-convert it into
    # a date serial number, matching the other rows in the column, then
    # restore the default (Normal) style so no stray formatting is left
    # on the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "22:59:02"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.5
    # Column G (Exit Price) stays blank, matching the other open trades.
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Column P (Exit Reason) stays blank, matching the other open trades.
    $ws.Cells.Item($row, 17).Value = 0
}
